# Added new Excel sheet to consolidate all 2017 ETNP MS/MS numbers
# -> rename the second sheet ("prok only") to "DB peps", and update the
#    saved view state of that sheet (scroll position / selection).

$wb = $excel.ActiveWorkbook

# Rename the "prok only" sheet to "DB peps"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "DB peps"

# Reset the view on the renamed sheet: scroll back to the top (clearing the
# previous topLeftCell="A378") and select cell D18 instead of the old
# A340:XFD340 row selection.
$ws2.Activate()
$ws2.Range("D18").Select()
